$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 41480
$ws.Range("J120").Value = 41480
$ws.Range("L120").Value = 41480
$ws.Range("N120").Value = -51156

$ws.Range("H132").Value = 1425983.2
$ws.Range("I132").Value = 1262.9701
$ws.Range("J132").Value = 10103825
$ws.Range("K132").Value = 3788.9103
$ws.Range("L132").Value = 30311475
$ws.Range("M132").Value = -1258.9103
$ws.Range("N132").Value = -30316535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 999.5
$ws.Range("I4").Value = 999
$ws.Range("K4").Value = 999
$ws.Range("M4").Value = -883

$ws.Range("H32").Value = 5006518.5
$ws.Range("I32").Value = 6811.3457
$ws.Range("J32").Value = 26321058
$ws.Range("K32").Value = 6811.3457
$ws.Range("L32").Value = 26321058
$ws.Range("M32").Value = -6524.3457
$ws.Range("N32").Value = -26321632

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 520.8
$ws.Range("I12").Value = 520.8
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 520.8
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -352.8
$ws.Range("N12").ClearContents()

$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H24").Value = 570
$ws.Range("I24").Value = 462.5
$ws.Range("K24").Value = 462.5
$ws.Range("M24").Value = -227.5

$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10228

$ws.Range("H44").Value = 5750
$ws.Range("I44").Value = 4000
$ws.Range("J44").Value = 6333.3335
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 6333.3335
$ws.Range("M44").Value = -3503
$ws.Range("N44").Value = -7327.3335

$ws.Range("H134").Value = 31513848
$ws.Range("I134").Value = 38462884
$ws.Range("J134").Value = 8929471
$ws.Range("K134").Value = 115388652
$ws.Range("L134").Value = 26788413
$ws.Range("M134").Value = -115386117
$ws.Range("N134").Value = -26793483

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 256.9375
$ws.Range("I5").Value = 90.111115
$ws.Range("K5").Value = 90.111115
$ws.Range("M5").Value = 21.888885

$ws.Range("H13").Value = 9486.429
$ws.Range("J13").Value = 5601.25
$ws.Range("L13").Value = 5601.25
$ws.Range("N13").Value = -5879.25

$ws.Range("H41").Value = 8332.5
$ws.Range("J41").Value = 8331.666999999999
$ws.Range("L41").Value = 8331.666999999999
$ws.Range("N41").Value = -9187.666999999999

$ws.Range("H42").Value = 10933.333

$ws.Range("H45").Value = 9000
$ws.Range("I45").Value = 9000
$ws.Range("J45").Value = 9000
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = -8407
$ws.Range("N45").Value = -10186

$ws.Range("H58").Value = 1568280.2
$ws.Range("I58").Value = 781.25
$ws.Range("J58").Value = 5051611.5
$ws.Range("K58").Value = 781.25
$ws.Range("L58").Value = 5051611.5
$ws.Range("M58").Value = -578.25
$ws.Range("N58").Value = -5052017.5

$ws.Range("H136").Value = 1568280.2
$ws.Range("I136").Value = 781.25
$ws.Range("J136").Value = 5051611.5
$ws.Range("K136").Value = 2343.75
$ws.Range("L136").Value = 15154834.5
$ws.Range("M136").Value = 206.25
$ws.Range("N136").Value = -15159934.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4241644
$ws.Range("I5").Value = 10989637
$ws.Range("J5").Value = 1755541.6
$ws.Range("K5").Value = 32968911
$ws.Range("L5").Value = 5266624.800000001
$ws.Range("M5").Value = -32968799
$ws.Range("N5").Value = -5266848.800000001

$ws.Range("H86").Value = 838.5714
$ws.Range("J86").Value = 857.1111
$ws.Range("L86").Value = 2571.3333
$ws.Range("N86").Value = -4943.3333

$ws.Range("H89").Value = 838.5714
$ws.Range("J89").Value = 857.1111
$ws.Range("L89").Value = 7713.9999
$ws.Range("N89").Value = -19569.9999

$ws.Range("H107").Value = 1221910.1
$ws.Range("I107").Value = 2331105
$ws.Range("J107").Value = 1795.8
$ws.Range("K107").Value = 6993315
$ws.Range("L107").Value = 5387.4
$ws.Range("M107").Value = -6991395
$ws.Range("N107").Value = -9227.4

$ws.Range("H135").Value = 4241644
$ws.Range("I135").Value = 10989637
$ws.Range("J135").Value = 1755541.6
$ws.Range("K135").Value = 98906733
$ws.Range("L135").Value = 15799874.4
$ws.Range("M135").Value = -98904198
$ws.Range("N135").Value = -15804944.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25770.75
$ws.Range("J57").Value = 25770.75
$ws.Range("L57").Value = 25770.75
$ws.Range("N57").Value = -27410.75

$ws.Range("H107").Value = 126.375
$ws.Range("I107").Value = 101.46667
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 101.46667
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1818.53333
$ws.Range("N107").Value = -4340

$ws.Range("H113").Value = 18907.072
$ws.Range("I113").Value = 790.8182
$ws.Range("J113").Value = 85333.336
$ws.Range("K113").Value = 790.8182
$ws.Range("L113").Value = 85333.336
$ws.Range("M113").Value = 1379.1818
$ws.Range("N113").Value = -89673.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 265
$ws.Range("I35").Value = 265
$ws.Range("K35").Value = 265
$ws.Range("M35").Value = 71

$ws.Range("H51").Value = 15730
$ws.Range("J51").Value = 15730
$ws.Range("L51").Value = 15730
$ws.Range("N51").Value = -16686

$ws.Range("H136").Value = 6537746
$ws.Range("I136").Value = 6537746
$ws.Range("K136").Value = 19613238
$ws.Range("M136").Value = -19610688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 4905
$ws.Range("J22").Value = 4905
$ws.Range("L22").Value = 4905
$ws.Range("N22").Value = -5491

$ws.Range("H34").Value = 8253.333000000001
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 8253.333000000001
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 8253.333000000001
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8659.333000000001

$ws.Range("H113").Value = 295.27777
$ws.Range("I113").Value = 230.66667
$ws.Range("J113").Value = 359.8889
$ws.Range("K113").Value = 692.00001
$ws.Range("L113").Value = 1079.6667
$ws.Range("M113").Value = 1477.99999
$ws.Range("N113").Value = -5419.6667
